# Weekly update: insert a new price record row for "Haba" (Femacal de La
# Calera market) right after the existing row 69, shifting all the
# subsequent rows (previously 70-86) down by one (now 71-87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70 (pushes old row 70.. down to 71..)
$ws.Rows("70:70").Insert()

# Populate the newly inserted row 70 with the new weekly record
$ws.Range("A70").Value2 = 3
$ws.Range("B70").Value2 = "Femacal de La Calera"
$ws.Range("C70").Value2 = "Coquimbo"
$ws.Range("D70").Value2 = 44504
$ws.Range("E70").Value2 = 5
$ws.Range("F70").Value2 = 100112026
$ws.Range("G70").Value2 = "Haba"
$ws.Range("H70").Value2 = "Sin especificar"
$ws.Range("I70").Value2 = "Primera"
$ws.Range("J70").Value2 = 120
$ws.Range("K70").Value2 = 7500
$ws.Range("L70").Value2 = 8000
$ws.Range("M70").Value2 = 7750
$ws.Range("N70").Value2 = '$/saco 25 kilos'
$ws.Range("O70").Value2 = "Provincia de Quillota"
$ws.Range("P70").Value2 = 310
$ws.Range("Q70").Value2 = 25
$ws.Range("R70").Value2 = "Hortaliza"
